$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @(19, "categorization", 4, 2, 1, 163, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_x0u5z.png", 92, 78.16216216216216, 85.08108108108108, 37, 10, 10, 10)
    3 = @(19, "categorization", 4, 2, 2, 164, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_le8uf.png", 12.88888888888889, 9.222222222222221, 11.05555555555556, 36, 1, 1, 1)
    4 = @(19, "categorization", 4, 2, 3, 165, "kitchens", "living_rooms", "distractor", $null, "f", "stimuli/img_pbsj1.png", 73.88636363636364, 51.52272727272727, 62.70454545454545, 44, 6, 6, 6)
    5 = @(19, "categorization", 4, 2, 4, 166, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_qgbyn.png", 65.08108108108108, 40.10810810810811, 52.5945945945946, 37, 4, 4, 4)
    6 = @(19, "categorization", 4, 2, 5, 167, "kitchens", "kitchens", "distractor", $null, "f", "stimuli/img_uegbb.png", 78.80952380952381, 61.52380952380953, 70.16666666666667, 42, 8, 8, 8)
    7 = @(19, "categorization", 4, 2, 6, 168, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_oou46.png", 75.70270270270271, 54.86486486486486, 65.28378378378379, 37, 6, 6, 6)
    8 = @(19, "categorization", 4, 2, 7, 169, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_ybbmx.png", 55.24324324324324, 36.75675675675676, 46, 37, 3, 3, 3)
    9 = @(19, "categorization", 4, 2, 8, 170, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_h0hbk.png", 86.80952380952381, 69.19047619047619, 78, 42, 9, 9, 9)
    10 = @(19, "categorization", 4, 2, 9, 171, "kitchens", "kitchens", "distractor", $null, "f", "stimuli/img_2b8fp.png", 73.89189189189189, 51.45945945945946, 62.67567567567568, 37, 6, 6, 6)
    11 = @(19, "categorization", 4, 2, 10, 172, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_wyctg.png", 33.44736842105263, 11.39473684210526, 22.42105263157895, 38, 1, 1, 1)
    12 = @(19, "categorization", 4, 2, 11, 173, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_a9acb.png", 77.11428571428571, 58.42857142857143, 67.77142857142857, 35, 7, 7, 7)
    13 = @(19, "categorization", 4, 2, 12, 174, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_71mhq.png", 69.34210526315789, 47.02631578947368, 58.18421052631579, 38, 5, 5, 5)
    14 = @(19, "categorization", 4, 2, 13, 175, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_5m6x4.png", 80.23076923076923, 58.41025641025641, 69.32051282051282, 39, 7, 7, 7)
    15 = @(19, "categorization", 4, 2, 14, 176, "kitchens", "kitchens", "distractor", $null, "f", "stimuli/img_a220l.png", 79.45945945945945, 60.97297297297298, 70.21621621621621, 37, 8, 8, 8)
    16 = @(19, "categorization", 4, 2, 15, 177, "kitchens", "kitchens", "distractor", $null, "f", "stimuli/img_q577a.png", 81.26470588235294, 59.08823529411764, 70.17647058823529, 34, 8, 8, 8)
    17 = @(19, "categorization", 4, 2, 16, 178, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_bj2gr.png", 65.25, 44.8, 55.025, 40, 4, 4, 4)
    18 = @(19, "categorization", 4, 2, 17, 179, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_uxxo0.png", 71.74418604651163, 48.44186046511628, 60.09302325581395, 43, 5, 5, 5)
    19 = @(19, "categorization", 4, 2, 18, 180, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_t2ioc.png", 88.1891891891892, 74.05405405405405, 81.12162162162161, 37, 10, 10, 10)
    20 = @(19, "categorization", 4, 2, 19, 181, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_rvssl.png", 74.25, 54.33333333333334, 64.29166666666667, 36, 6, 6, 6)
    21 = @(19, "categorization", 4, 2, 20, 182, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_2js6m.png", 40.02777777777778, 20.88888888888889, 30.45833333333334, 36, 2, 2, 2)
    22 = @(19, "categorization", 4, 2, 21, 183, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_th7xh.png", 82.35897435897436, 65.53846153846153, 73.94871794871796, 39, 8, 8, 8)
    23 = @(19, "categorization", 4, 2, 22, 184, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_v8dra.png", 61.77272727272727, 38.79545454545455, 50.28409090909091, 44, 3, 3, 3)
    24 = @(19, "categorization", 4, 2, 23, 185, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_okvvw.png", 50.58333333333334, 32.11111111111111, 41.34722222222223, 36, 2, 2, 2)
    25 = @(19, "categorization", 4, 2, 24, 186, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_2pk6v.png", 85.08108108108108, 66.16216216216216, 75.62162162162161, 37, 9, 9, 9)
    26 = @(19, "categorization", 4, 2, 25, 187, "kitchens", "bedrooms", "target", $null, "j", "stimuli/img_fqgem.png", 80.75, 61.475, 71.1125, 40, 8, 8, 8)
    27 = @(19, "categorization", 4, 2, 26, 188, "kitchens", "kitchens", "distractor", $null, "f", "stimuli/img_cv6mf.png", 66.8, 42.08, 54.44, 25, 4, 4, 4)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        if ($null -ne $vals[$c]) {
            $ws.Cells.Item([int]$r, $c + 1).Value = $vals[$c]
        }
    }
}

# Row 28 (subject_id=19, trial_total=270) was removed in the edit.
$ws.Rows(28).Delete()
